$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old leading key column (A); B:F shift left to become A:E
$ws.Range("A1").EntireColumn.Delete()

# Fix header text typo, now in column D: MODEL_CONDITION -> MODELCONDITION
$ws.Range("D1").Value = "MODELCONDITION"
